$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the value to be stored as text (not auto-converted to a number
    # or losing formatting such as trailing zeros), while preserving the
    # cell's original style (avoids introducing a "Text" number format).
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '65.860.96'
Set-TextValue $ws.Range("E2") '  +1.12%  '
Set-TextValue $ws.Range("D3") '2.694.75'
Set-TextValue $ws.Range("E3") '  +2.00%  '
Set-TextValue $ws.Range("E4") '  -0.10%  '
Set-TextValue $ws.Range("D5") '607.41'
Set-TextValue $ws.Range("E5") '  +1.75%  '
Set-TextValue $ws.Range("D6") '157.81'
Set-TextValue $ws.Range("E6") '  +1.33%  '
Set-TextValue $ws.Range("E7") '  -0.05%  '
Set-TextValue $ws.Range("E8") '  -0.84%  '
Set-TextValue $ws.Range("E9") '  +5.22%  '
Set-TextValue $ws.Range("D10") '6.00'
Set-TextValue $ws.Range("E10") '  +3.53%  '
Set-TextValue $ws.Range("D11") '0.402'
Set-TextValue $ws.Range("E11") '  -0.13%  '
Set-TextValue $ws.Range("E12") '  +0.94%  '
Set-TextValue $ws.Range("D13") '30.28'
Set-TextValue $ws.Range("E13") '  +3.69%  '
Set-TextValue $ws.Range("E14") '  +8.03%  '
Set-TextValue $ws.Range("D15") '3.178.73'
Set-TextValue $ws.Range("E15") '  +1.95%  '
Set-TextValue $ws.Range("D16") '65.699.05'
Set-TextValue $ws.Range("E16") '  +1.01%  '
Set-TextValue $ws.Range("D17") '2.683.40'
Set-TextValue $ws.Range("E17") '  +1.75%  '
Set-TextValue $ws.Range("D18") '12.66'
Set-TextValue $ws.Range("E18") '  +0.59%  '
Set-TextValue $ws.Range("E19") '  +1.98%  '
Set-TextValue $ws.Range("D20") '358.37'
Set-TextValue $ws.Range("E20") '  +1.33%  '
Set-TextValue $ws.Range("D21") '7.54'
Set-TextValue $ws.Range("E21") '  +2.45%  '
Set-TextValue $ws.Range("E22") '  -0.20%  '
Set-TextValue $ws.Range("D23") '70.60'
Set-TextValue $ws.Range("E23") '  +3.50%  '
Set-TextValue $ws.Range("D24") '9.83'
Set-TextValue $ws.Range("E24") '  +2.79%  '
Set-TextValue $ws.Range("B25") 'PEPE'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D25") '0.0000106'
Set-TextValue $ws.Range("E25") '  +12.14%  '
Set-TextValue $ws.Range("B26") 'SuiNetwork'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D26") '1.66'
Set-TextValue $ws.Range("E26") '  -2.37%  '
Set-TextValue $ws.Range("E27") '  +2.36%  '
Set-TextValue $ws.Range("D28") '0.171'
Set-TextValue $ws.Range("E28") '  +3.97%  '
Set-TextValue $ws.Range("D29") '8.39'
Set-TextValue $ws.Range("E29") '  +3.63%  '
Set-TextValue $ws.Range("E30") '  +4.56%  '
Set-TextValue $ws.Range("B31") 'Bittensor'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D31") '540.37'
Set-TextValue $ws.Range("E31") '  +6.60%  '
Set-TextValue $ws.Range("B32") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D32") '1.00'
Set-TextValue $ws.Range("E32") '  +0.06%  '
Set-TextValue $ws.Range("E33") '  +1.81%  '
Set-TextValue $ws.Range("D34") '6.71'
Set-TextValue $ws.Range("E34") '  +5.66%  '
Set-TextValue $ws.Range("D35") '5.41'
Set-TextValue $ws.Range("E35") '  -3.21%  '
Set-TextValue $ws.Range("E36") '  +1.90%  '
Set-TextValue $ws.Range("D37") '20.84'
Set-TextValue $ws.Range("E37") '  +2.65%  '
Set-TextValue $ws.Range("D38") '163.57'
Set-TextValue $ws.Range("E38") '  -0.08%  '
Set-TextValue $ws.Range("E39") '  -0.48%  '
Set-TextValue $ws.Range("E40") '  -0.13%  '
Set-TextValue $ws.Range("D41") '171.34'
Set-TextValue $ws.Range("E41") '  +3.19%  '
Set-TextValue $ws.Range("E42") '  +0.00%  '
Set-TextValue $ws.Range("E43") '  +0.32%  '
Set-TextValue $ws.Range("D44") '4.17'
Set-TextValue $ws.Range("E44") '  +1.97%  '
Set-TextValue $ws.Range("E45") '  -0.10%  '
Set-TextValue $ws.Range("D46") '23.51'
Set-TextValue $ws.Range("E46") '  +1.93%  '
Set-TextValue $ws.Range("D47") '2.27'
Set-TextValue $ws.Range("E47") '  +4.02%  '
Set-TextValue $ws.Range("D48") '0.0266'
Set-TextValue $ws.Range("E48") '  +3.93%  '
Set-TextValue $ws.Range("D49") '0.654'
Set-TextValue $ws.Range("E49") '  +0.96%  '
Set-TextValue $ws.Range("D50") '20.82'
Set-TextValue $ws.Range("E50") '  +6.79%  '
Set-TextValue $ws.Range("D51") '0.0991'
Set-TextValue $ws.Range("E51") '  +0.65%  '
